# Update "想去人数" (interest count) and "最低票价" (min price) figures
# on the "展览" and "全部类型" worksheets to match the refreshed scrape
# output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 233
    $ws.Range("G4").Value = 55
    $ws.Range("F7").Value = 5753
    $ws.Range("F8").Value = 5095
}

# Row with the "MAX特摄同人only2.0" event lives at row 14 on "展览"
# and row 16 on "全部类型".
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F14").Value = 18

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F16").Value = 18
